$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 487, pushing existing rows 487:553 down to 488:554.
$ws.Rows.Item(487).Insert()

# Populate the newly inserted row 487 with the new data entry.
$ws.Cells.Item(487, 1).Value = 3
$ws.Cells.Item(487, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(487, 3).Value = "Coquimbo"
$ws.Cells.Item(487, 4).Value = 45077
$ws.Cells.Item(487, 5).Value = 5
$ws.Cells.Item(487, 6).Value = 100112009
$ws.Cells.Item(487, 7).Value = "Acelga"
$ws.Cells.Item(487, 8).Value = "Sin especificar"
$ws.Cells.Item(487, 9).Value = "Primera"
$ws.Cells.Item(487, 10).Value = 230
$ws.Cells.Item(487, 11).Value = 3300
$ws.Cells.Item(487, 12).Value = 3500
$ws.Cells.Item(487, 13).Value = 3404
$ws.Cells.Item(487, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(487, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(487, 16).Value = 567
$ws.Cells.Item(487, 17).Value = 6
$ws.Cells.Item(487, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Cells.Item(487, 4).NumberFormat = $ws.Cells.Item(488, 4).NumberFormat
